# Add ability for the user to set "days off" -- entries whose Date matches
# one of the days-off values are skipped when the schedule is written out,
# and everything after them shifts up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day(s) the user marked as missed/off. (Here: 4/10.)
$daysOff = @("4/10")

$firstRow = 4
$lastRow = 50

# 1. Read every existing entry back into one ordered stream: for each row,
#    the left block (A,B,C) first, then the right block (E,F,G).
$entries = New-Object System.Collections.ArrayList

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dateVal = $ws.Range("A" + $r).Value()
    if ($dateVal -ne $null -and $dateVal -ne "") {
        $entry = @{
            Date   = $dateVal
            Class  = $ws.Range("B" + $r).Value()
            Length = $ws.Range("C" + $r).Value()
        }
        [void]$entries.Add($entry)
    }

    $dateVal2 = $ws.Range("E" + $r).Value()
    if ($dateVal2 -ne $null -and $dateVal2 -ne "") {
        $entry2 = @{
            Date   = $dateVal2
            Class  = $ws.Range("F" + $r).Value()
            Length = $ws.Range("G" + $r).Value()
        }
        [void]$entries.Add($entry2)
    }
}

# 2. Drop any entry whose date is a day the user marked off.
$kept = New-Object System.Collections.ArrayList
foreach ($entry in $entries) {
    if (-not ($daysOff -contains $entry.Date)) {
        [void]$kept.Add($entry)
    }
}

# 3. Clear the old data area (rows firstRow..lastRow, columns A:C and E:G).
$clearRange = "A" + $firstRow + ":G" + $lastRow
$ws.Range($clearRange).Clear()

# 4. Write the kept entries back, two per row (left block then right block),
#    starting at firstRow.
$row = $firstRow
$i = 0
$count = $kept.Count
while ($i -lt $count) {
    $left = $kept[$i]
    $ws.Range("A" + $row).Value = $left.Date
    $ws.Range("B" + $row).Value = $left.Class
    $ws.Range("C" + $row).Value = $left.Length
    $i = $i + 1

    if ($i -lt $count) {
        $right = $kept[$i]
        $ws.Range("E" + $row).Value = $right.Date
        $ws.Range("F" + $row).Value = $right.Class
        $ws.Range("G" + $row).Value = $right.Length
        $i = $i + 1
    }

    $row = $row + 1
}

$lastUsedRow = $row - 1

# 5. Clear any now-unused trailing rows (the grid shrank).
if ($lastUsedRow -lt $lastRow) {
    $ws.Range("A" + ($lastUsedRow + 1) + ":H" + $lastRow).Clear()
}
